$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "288.42"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "1.05%"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "29.21"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "1.97%"
$r.Style = "Normal"
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "5.079"
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "3.31%"
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.06680"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "2.92%"
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "7.332"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "1.56%"
$r.Style = "Normal"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "3.405"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "1.10%"
$r.Style = "Normal"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "1.352"
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "1.21%"
$r.Style = "Normal"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.9184"
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "0.50%"
$r.Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.1588"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "3.07%"
$r.Style = "Normal"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.06816"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "7.62%"
$r.Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.07669"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "0.40%"
$r.Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.02935"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "-1.36%"
$r.Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.08988"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "0.28%"
$r.Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.001587"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "-1.24%"
$r.Style = "Normal"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.04497"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "0.61%"
$r.Style = "Normal"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.0006473"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "-1.09%"
$r.Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.006248"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "3.69%"
$r.Style = "Normal"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "3.452"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "-0.26%"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "-1.07%"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "2.03%"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "-2.46%"
$r.Style = "Normal"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "4.056"
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "1.50%"
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.1581"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "1.62%"
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.001191"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "0.13%"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.004116"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "-4.92%"
$r.Style = "Normal"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "1.59%"
$r.Style = "Normal"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.0001616"
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "-1.20%"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.04230"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "1.76%"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.006719"
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "0.54%"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.1238"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "0.49%"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "-3.69%"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "13.62%"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.00005714"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "6.36%"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.968"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "-3.59%"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.01306"
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "-29.45%"
$r.Style = "Normal"
